$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, country name (col A), then Casos totales,
# Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes (cols B-H)
$updates = @(
    @(4, 'Estados Unidos', 6089000, 42366, 3367910, 2535377, 0, 917, 185713),
    @(5, 'Brasil', 3812605, 48112, 2947250, 745784, 0, 845, 119571),
    @(9, 'Sudafrica', 620132, 1846, 533935, 72454, 0, 115, 13743),
    @(34, 'Egipto', 98285, 223, 70419, 22504, 0, 20, 5362),
    @(79, 'Costa de Marfil', 17797, 95, 16315, 1367, 0, 0, 115),
    @(81, 'Bulgaria', 16065, 157, 11231, 4231, 0, 9, 603),
    @(120, 'Ruanda', 3742, 70, 1866, 1860, 0, 1, 16),
    @(121, 'Eslovaquia', 3728, 102, 2225, 1470, 0, 0, 33),
    @(155, 'Togo', 1365, 39, 981, 357, 0, 0, 27),
    @(156, 'Burkina Faso', 1352, 0, 1058, 239, 0, 0, 55),
    @(158, 'Guyana', 1180, 40, 633, 512, 0, 3, 35),
    @(159, 'Niger', 1173, 0, 1084, 20, 0, 0, 69),
    @(160, 'Guadalupe', 1145, 210, 336, 793, 0, 1, 16),
    @(161, 'Principado de Andorra', 1124, 26, 902, 169, 0, 0, 53),
    @(162, 'Lesoto', 1051, 0, 526, 494, 0, 0, 31),
    @(163, 'Vietnam', 1038, 2, 663, 345, 0, 0, 30),
    @(164, 'Republica del Chad', 1008, 4, 878, 53, 0, 0, 77),
    @(166, 'Belice', 870, 52, 93, 765, 0, 0, 12),
    @(186, 'San Martin (Parte Francesa)', 213, 15, 79, 129, 0, 0, 5),
    @(187, 'Islas Caimanes', 205, 0, 202, 2, 0, 0, 1),
    @(190, 'Barbados', 166, 1, 141, 18, 0, 0, 7),
    @(197, 'Curazao', 57, 4, 35, 21, 0, 0, 1),
    @(209, 'San Bartolome', 18, 1, 13, 5, 0, 0, 0),
    @(210, 'Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0)
)

$cols = @("B","C","D","E","F","G","H")

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Range("A" + $r).Value = $u[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $u[2 + $i]
    }
}

# Update the "last refreshed" timestamp shown in A1
$ws.Range("A1").Value = 'Datos actualizados a 28 de Agosto de 2020 a las 23:59'
